$d = $word.ActiveDocument

# The target text currently lives in a single run. We need to split it into
# two runs (same formatting) without leaving any tracked-change markup
# behind. Word's Find/Replace and Range.Text/InsertAfter normally just
# rewrite text within a run's existing boundaries (merging same-formatted
# text), so we instead perform the edit under Track Changes -- which forces
# a real run boundary at the edit point -- and then Accept just the
# revisions we created (not Document.AcceptAllRevisions, which would touch
# unrelated parts of the document).

$oldTail = "has been replaced with a simple high score tracker"
$newTail = "has been implemented and am working towards adding more obstacles before TP3"

$d.TrackRevisions = $true

# Insert the new tail text immediately before the old tail text.
$rng = $d.Content
$found = $rng.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence tail to replace"
}
$rng.InsertBefore($newTail)

# Re-find (position shifted) and delete the old tail text.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not re-find old sentence tail to delete"
}
$rng2.Delete()

$d.TrackRevisions = $false

# Accept only the revisions we just made, one by one, so the rest of the
# document (rsids, lastRenderedPageBreak hints, etc.) is left untouched.
for ($i = $d.Revisions.Count; $i -ge 1; $i--) {
    $d.Revisions.Item($i).Accept()
}
